$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the format of A2 (border/bold/centered style) down to the new A3 and A4 cells
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Set column A group-name labels first (matches shared-string creation order)
$ws.Range("A2").Value = 'Agua y Saneamiento'
$ws.Range("A3").Value = 'INVEMAR - Calidad Ambiental Marina'
$ws.Range("A4").Value = 'Limnología y Recursos Hídricos'

# Then set column B long descriptions
$ws.Range("B2").Value = '7. Consultoría científica: Formulacion del Plan de Manejo del acuifero de la cuenca del rio Risaralda, implementación de los planes de manejo de los acuiferos Santagueda Km 41- Irra y rio grande de la Magdalena y operación de la red de monitoreo de agua Subterranea de los acuiferos en jurisdicción de Caldas  Año de inicio: 2015, Mes de inicio: 8, Año de fin: 2016, Mes de fin: 7  Idioma: Español, Ciudad: MANIZALES, Disponibilidad: No restringido, Duración: 0, Número del contrato: 136-2015, Institución en la cual prestó el servicio: Corporación Autónoma Regional de Caldas 
 22. Consultoría científica: Caracterización de agua subterránea para once (11) pozos de monitoreo ubicados en el Departamento de Risaralda. 2019-I  Año de inicio: 2019, Mes de inicio: 2, Año de fin: 2019, Mes de fin: 5  Idioma: Español, Ciudad: LA VIRGINIA, Disponibilidad: Restringido, Duración: 0, Número del contrato: 024-19, Institución en la cual prestó el servicio: INGENIO RISARALDA S.A. 
 26. Consultoría científica: Caracterización de agua subterránea para once (11) pozos de monitoreo ubicados en el Departamento de Risaralda. 2019-II  Año de inicio: 2019, Mes de inicio: 9, Año de fin: 2019, Mes de fin: 11  Idioma: Español, Ciudad: PEREIRA, Disponibilidad: Restringido, Duración: 0, Número del contrato: 103-19, Institución en la cual prestó el servicio: INGENIO RISARALDA S.A. 
 30. Consultoría científica: Caracterización de agua subterránea proveniente de 11 pozos de monitoreo ubicados en el Departamento de Risaralda. 2020  Año de inicio: 2020, Mes de inicio: 3, Año de fin: 2020, Mes de fin: 9  Idioma: Español, Ciudad: LA VIRGINIA, Disponibilidad: Restringido, Duración: 0, Número del contrato: 022-20, Institución en la cual prestó el servicio: INGENIO RISARALDA S.A. 
 34. Consultoría científica: Estudios de caracterización de aguas de pozos de monitoreo año 2016  Año de inicio: 2016, Mes de inicio: 2, Año de fin: 2016, Mes de fin: 5  Idioma: Español, Ciudad: LA VIRGINIA, Disponibilidad: Restringido, Duración: 0, Número del contrato: 018-16, Institución en la cual prestó el servicio: INGENIO RISARALDA S.A. 
 35. Consultoría científica: Estudio de caracterización de agua de pozos de monitoreo en el Ingenio Risaralda S.A, segundo semestre de 2016  Año de inicio: 2016, Mes de inicio: 2, Año de fin: 2016, Mes de fin: 9  Idioma: Español, Ciudad: LA VIRGINIA, Disponibilidad: Restringido, Duración: 7, Número del contrato: 106-16, Institución en la cual prestó el servicio: INGENIO RISARALDA S.A. 
 43. Consultoría científica: Estudio de caracterización de agua de pozos de monitoreo semestre I año 2018 Ingenio Risaralda S.A  Año de inicio: 2018, Mes de inicio: 3, Año de fin: 2018, Mes de fin: 6  Idioma: Español, Ciudad: LA VIRGINIA, Disponibilidad: Restringido, Duración: 0, Número del contrato: Orden de compra # 3620172063, Institución en la cual prestó el servicio: INGENIO RISARALDA S.A. 
 45. Consultoría científica: Estudio de caracterización de agua de pozos de monitoreo semestre II año 2018 Ingenio Risaralda S.A.  Año de inicio: 2018, Mes de inicio: 8, Año de fin: 2018, Mes de fin: 12  Idioma: Español, Ciudad: LA VIRGINIA, Disponibilidad: Restringido, Duración: 0, Número del contrato: Orden de compra # 3620178418, Institución en la cual prestó el servicio: INGENIO RISARALDA S.A. 
 51. Consultoría científica: Formulación del Plan de Manejo del acuífero de la cuenca del rio Risaralda, implementación de los planes de manejo de los acuíferos Santagueda Km 41- Irra y rio grande de la Magdalena y operación de la red de monitoreo de agua Subterránea de los acuíferos en jurisdicción de Caldas  Año de inicio: 2015, Mes de inicio: 8, Año de fin: 2016, Mes de fin: 9  Idioma: Español, Ciudad: , Disponibilidad: No restringido, Duración: 0, Número del contrato: Contrato 136-2015, Institución en la cual prestó el servicio: Corporación Autónoma Regional de Caldas'
$ws.Range("B3").Value = '5. Consultoría científica: Monitoreo de calidad de aguas, sedimentos y aspectos biológicos en una zona industrial de Mamonal, bahía de Cartagena  Año de inicio: 2007, Mes de inicio: 7, Año de fin: 2008, Mes de fin: 7  Idioma: Español, Ciudad: SANTA MARTA, Disponibilidad: Restringido, Duración: 12, Número del contrato: CONTRATO DE CONSULTORÍA - 2007, Institución en la cual prestó el servicio: INSTITUTO DE INVESTIGACIONES MARINAS Y COSTERAS JOSE BENITO VIVES DE ANDREIS INVEMAR 
 6. Consultoría científica: MONITOREO AMBIENTAL DE LA CALIDAD DE AGUAS, SEDIMENTOS Y COMUNIDADES MARINAS EN LA ZONA DE INFLUENCIA DIRECTA DEL POZO EXPLORATORIO ARAZÁ I  Año de inicio: 2007, Mes de inicio: 6, Año de fin: 2008, Mes de fin: 7  Idioma: Español, Ciudad: SANTA MARTA, Disponibilidad: Restringido, Duración: 12, Número del contrato: , Institución en la cual prestó el servicio: INSTITUTO DE INVESTIGACIONES MARINAS Y COSTERAS JOSE BENITO VIVES DE ANDREIS INVEMAR 
 13. Consultoría científica: Monitoreo de calidad de aguas marinas, potables y residuales durante la perforación del Pozo Brama-01  Año de inicio: 2017, Mes de inicio: 4, Año de fin: 2018, Mes de fin: 2  Idioma: Español, Ciudad: SANTA MARTA, Disponibilidad: Restringido, Duración: 10, Número del contrato: Contrato PETROBRAS -INVEMAR No. 9019234, Institución en la cual prestó el servicio: INSTITUTO DE INVESTIGACIONES MARINAS Y COSTERAS JOSE BENITO VIVES DE ANDREIS INVEMAR 
 30. Consultoría científica: Monitoreo de Calidad de Aguas, Sedimentos, Fauna Asociada a Manglar y Bentos en el Área de Influencia del Dragado en el Puerto de Buenaventura  Año de inicio: 2008, Mes de inicio: 1, Año de fin: 0, Mes de fin: 0  Idioma: Español, Ciudad: SANTA MARTA, Disponibilidad: Restringido, Duración: 18, Número del contrato: , Institución en la cual prestó el servicio: Ondenemeingen Jande Nul NV 
 35. Consultoría científica: PROTOCOLO DE MONITOREO DE LAS AGUAS DE LASTRE DE LOS BUQUES QUE INGRESAN AL PUERTO DE SANTA MARTA  Año de inicio: 2012, Mes de inicio: 9, Año de fin: 2014, Mes de fin: 2  Idioma: Español, Ciudad: SANTA MARTA, Disponibilidad: Restringido, Duración: 17, Número del contrato: , Institución en la cual prestó el servicio: INSTITUTO DE INVESTIGACIONES MARINAS Y COSTERAS JOSE BENITO VIVES DE ANDREIS INVEMAR'
$ws.Range("B4").Value = '3. Consultoría científica: Monitoreo y seguimiento hidrobiológico del trasvase del río Guarinó y Monitoreo limnológico e hidrobiológico embalse Amaní - río la Miel y afluentes principales en el sector aguas abajo de la presa de la central hidroeléctrica Miel I - 2013  Año de inicio: 2011, Mes de inicio: 1, Año de fin: 2013, Mes de fin: 1  Idioma: Español, Ciudad: MEDELLÍN, Disponibilidad: Restringido, Duración: 36, Número del contrato: 46-3643, Institución en la cual prestó el servicio: ISAGEN S.A. ESP 
 9. Consultoría científica: Monitoreo limnológico e hidrobiológico embalse Amaní - río la Miel y afluentes principales en el sector aguas abajo de la presa de la central hidroeléctrica Miel I - 2009 Año de inicio: 2008, Mes de inicio: 1, Año de fin: 2009, Mes de fin: 12  Idioma: Español, Ciudad: MEDELLÍN, Disponibilidad: Restringido, Duración: 24, Número del contrato: 46-2628, Institución en la cual prestó el servicio: ISAGEN S.A. ESP 
 22. Consultoría científica: Monitoreo Limnológico e hidrobiológico embalse Amaní-Río La Miel y afluentes principales en el sector aguas abajo de la presa de la Central Hidroeléctrica Miel I - PMA 2006  Año de inicio: 2006, Mes de inicio: 1, Año de fin: 2007, Mes de fin: 1  Idioma: Español, Ciudad: , Disponibilidad: Restringido, Duración: 12, Número del contrato: 1, Institución en la cual prestó el servicio: ISAGEN S.A. ESP 
 24. Consultoría científica: Programa de monitoreo limnológico embalses Punchiná y San Lorenzo, Centrales Hidroeléctricas San Carlos y Jaguas PMA 2007  Año de inicio: 2007, Mes de inicio: 2, Año de fin: 0, Mes de fin: 1  Idioma: Español, Ciudad: , Disponibilidad: Restringido, Duración: 11, Número del contrato: , Institución en la cual prestó el servicio: ISAGEN S.A. ESP 
 32. Consultoría científica: Monitoreo Limnológico e hidrobiológico embalse Amaní-Río La Miel y afluentes principales en el sector aguas abajo de la presa de la Central Hidroeléctrica Miel I - 2012  Año de inicio: 2010, Mes de inicio: 1, Año de fin: 2012, Mes de fin: 1  Idioma: Español, Ciudad: RIONEGRO, Disponibilidad: Restringido, Duración: 24, Número del contrato: 46-3346, Institución en la cual prestó el servicio: ISAGEN S.A. ESP 
 46. Consultoría científica: Monitoreo de ictiofauna del río la miel y afluentes principales en el sector aguas abajo de la presa de la central hidroeléctrica miel I  Año de inicio: 2010, Mes de inicio: 1, Año de fin: 0, Mes de fin: 1  Idioma: Español, Ciudad: NORCASIA, Disponibilidad: Restringido, Duración: 0, Número del contrato: 46/3346, Institución en la cual prestó el servicio: UNIVERSIDAD CATOLICA DE ORIENTE 
 47. Consultoría científica: Monitoreo de ictiofauna del río la miel y afluentes principales en el sector aguas abajo de la presa de la central hidroeléctrica miel I-2008  Año de inicio: 2008, Mes de inicio: 1, Año de fin: 0, Mes de fin: 1  Idioma: Español, Ciudad: NORCASIA, Disponibilidad: Restringido, Duración: 0, Número del contrato: 46/2628, Institución en la cual prestó el servicio: UNIVERSIDAD CATOLICA DE ORIENTE'

# The long multi-line descriptions force Excel to auto-expand row height;
# reset rows back to the default/standard height so no explicit row height is stored
$ws.Rows(2).EntireRow.AutoFit()
$ws.Rows(3).EntireRow.AutoFit()
$ws.Rows(4).EntireRow.AutoFit()
